# caso 4 restriccion comentada x2
# Update the "Resultados" sheet with the new solver output values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resultados")

# Row 4 (interval I1)
$ws.Range("D4").Value = 12
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0
$ws.Range("I4").Value = 126600

# Row 5 (interval I2)
$ws.Range("D5").Value = 12
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 0

# Row 6 (interval I3)
$ws.Range("D6").Value = 12
$ws.Range("E6").Value = 8
$ws.Range("F6").Value = 0
$ws.Range("I6").Value = 237600

# Row 7 (interval I4)
$ws.Range("D7").Value = 12
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 2

# Row 8 (interval I5)
$ws.Range("D8").Value = 12
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = 0
$ws.Range("I8").Value = 260100

# Minimum MW price updates (S4 block)
$ws.Range("G11").Value = 8.44
$ws.Range("G13").Value = 9.504
$ws.Range("G15").Value = 9.633333333333333
